$d = $word.ActiveDocument

$replacements = @(
    @("55×89=", "33×47="),
    @("84×25=", "51×91="),
    @("26×91=", "48×67="),
    @("83×39=", "93×97="),
    @("43×56=", "52×12="),
    @("53×62=", "62×67="),
    @("14×47=", "66×72="),
    @("79×63=", "66×96="),
    @("31×42=", "29×54="),
    @("55×21=", "64×46="),
    @("71×13=", "93×24="),
    @("36×72=", "19×74="),
    @("80×78=", "61×11="),
    @("38×68=", "49×57="),
    @("54×32=", "35×44="),
    @("80×88=", "58×63="),
    @("88×58=", "16×24="),
    @("54×68=", "57×61="),
    @("11×90=", "50×94="),
    @("59×56=", "37×44="),
    @("87×83=", "40×26="),
    @("39×74=", "89×95="),
    @("27×31=", "31×27="),
    @("82×41=", "28×97="),
    @("61×68=", "74×76=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
